$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 2 updates ---
$ws.Range("B2").Value = "Binance"
$ws.Range("D2").Value = 44501
$ws.Range("E2").Value = 44561

# --- Row 3 updates ---
$ws.Range("C3").Value = "BTCUSDT"
$ws.Range("D3").Value = 44501
$ws.Range("E3").Value = 44561

# --- Rows 4 and 5: clear back down to the "blank template" row pattern ---
foreach ($r in 4,5) {
    $ws.Range("A$r").Clear()
    $ws.Range("H$r`:I$r").Clear()

    $ws.Range("B$r").ClearContents()
    $ws.Range("C$r").ClearContents()
    $ws.Range("D$r").ClearContents()
    $ws.Range("E$r").ClearContents()
    $ws.Range("F$r").ClearContents()
    $ws.Range("G$r").ClearContents()
    $ws.Range("J$r").ClearContents()

    $ws.Range("B6:G6").Copy()
    $ws.Range("B$r`:G$r").PasteSpecial(-4122)
    $ws.Range("J6").Copy()
    $ws.Range("J$r").PasteSpecial(-4122)
}

# --- Remove the trailing blank rows 126:128 ---
$ws.Rows("126:128").Delete()

# --- Update selection ---
$ws.Range("B3").Select()
